# PreCalc_Day_070 Test.pptx — apply commit "update rational exponents, add limits notes"
#
# Changes:
#  1. Date placeholder fields (handout master, slide master, all 11 slide
#     layouts) refresh their cached "datetimeFigureOut" text from 4/10/2018
#     to 4/11/2018.
#  2. Slide 2 ("Bell Work" title): adjacent runs with identical formatting
#     get merged by PowerPoint's save pass (no visible text change) --
#     "Bell Work: " + "#1-2 " -> "Bell Work: #1-2 ", and the six runs after
#     the line break collapse into a single run.
#  3. Slide 3 title: "Review Assignment (do what you need)" is edited down
#     to "Review Assignment" (split across two runs).

$p = $ppt.ActivePresentation

function Set-DateFieldText {
    param($shape, [string]$newText)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -ne $newText) {
            $tr.Characters(1, $tr.Length).Text = $newText
        }
    }
}

# --- 1. Refresh the cached date text everywhere the date placeholder lives ---

# Handout master
for ($i = 1; $i -le $p.HandoutMaster.Shapes.Count; $i++) {
    $shp = $p.HandoutMaster.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        Set-DateFieldText $shp "4/11/2018"
    }
}

# Slide master
for ($i = 1; $i -le $p.SlideMaster.Shapes.Count; $i++) {
    $shp = $p.SlideMaster.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        Set-DateFieldText $shp "4/11/2018"
    }
}

# Every slide layout
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*") {
            Set-DateFieldText $shp "4/11/2018"
        }
    }
}

# --- 2. Slide 2: "Bell Work" title textbox -> merge runs (text unchanged) ---

$slide2 = $p.Slides.Item(2)
$bellWork = $slide2.Shapes.Item(1)
$tr = $bellWork.TextFrame.TextRange

# "Bell Work: " + "#1-2 " -> single run "Bell Work: #1-2 "
$tr.Characters(1, 16).Text = "Bell Work: #1-2 "

# The six runs after the manual line break collapse into one run.
$tail = "	   	     #3-4 Solve    [L3]	#5 Solve[L4]"
$tr.Characters(31, 41).Text = $tail

# --- 3. Slide 3: title text trimmed to "Review Assignment" ---

$slide3 = $p.Slides.Item(3)
$reviewTitle = $slide3.Shapes.Item(1)
$tr3 = $reviewTitle.TextFrame.TextRange
$tr3.Characters(1, $tr3.Length).Text = "Review "
$tr3.InsertAfter("Assignment")
